$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header row (English labels), preserving the column order so the
# shared-string table indices line up with the target file:
#   A1 -> "Problem Component Name"
#   C1 -> "Status"
#   B1 -> "Note"
$ws.Range("A1").Value = "Problem Component Name"
$ws.Range("C1").Value = "Status"
$ws.Range("B1").Value = "Note"

# Move the active selection to B2 (matches the saved workbook state)
$ws.Range("B2").Select()
